$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.035.77'
$ws.Range("E2").Value = '  +7.83%  '
$ws.Range("D3").Value = '3.512.17'
$ws.Range("E3").Value = '  +11.16%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '190.30'
$ws.Range("E5").Value = '  +12.31%  '
$ws.Range("D6").Value = '550.93'
$ws.Range("E6").Value = '  +5.98%  '
$ws.Range("D7").Value = '3.504.75'
$ws.Range("E7").Value = '  +11.01%  '
$ws.Range("E8").Value = '  +3.33%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '0.632'
$ws.Range("E10").Value = '  +5.97%  '
$ws.Range("E11").Value = '  +17.68%  '
$ws.Range("D12").Value = '54.97'
$ws.Range("E12").Value = '  +5.76%  '
$ws.Range("E13").Value = '  +8.82%  '
$ws.Range("E14").Value = '  +4.41%  '
$ws.Range("D15").Value = '4.068.46'
$ws.Range("E15").Value = '  +11.42%  '
$ws.Range("D16").Value = '3.510.34'
$ws.Range("E16").Value = '  +11.52%  '
$ws.Range("E17").Value = '  +4.63%  '
$ws.Range("D18").Value = '67.004.59'
$ws.Range("E18").Value = '  +8.21%  '
$ws.Range("D19").Value = '18.17'
$ws.Range("E19").Value = '  +6.89%  '
$ws.Range("D20").Value = '11.88'
$ws.Range("E20").Value = '  +9.50%  '
$ws.Range("E21").Value = '  +3.41%  '
$ws.Range("D22").Value = '426.33'
$ws.Range("E22").Value = '  +18.57%  '
$ws.Range("E23").Value = '  +5.72%  '
$ws.Range("D24").Value = '84.67'
$ws.Range("E24").Value = '  +5.97%  '
$ws.Range("E25").Value = '  +7.10%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '2.90'
$ws.Range("E27").Value = '  +11.73%  '
$ws.Range("D28").Value = '11.97'
$ws.Range("E28").Value = '  +7.56%  '
$ws.Range("D29").Value = '8.93'
$ws.Range("E29").Value = '  +10.91%  '
$ws.Range("D30").Value = '30.16'
$ws.Range("E30").Value = '  +7.92%  '
$ws.Range("D31").Value = '647.83'
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("E32").Value = '  +5.02%  '
$ws.Range("E33").Value = '  +4.28%  '
$ws.Range("E34").Value = '  +6.55%  '
$ws.Range("D35").Value = '59.23'
$ws.Range("E35").Value = '  +5.40%  '
$ws.Range("D36").Value = '38.48'
$ws.Range("E36").Value = '  +5.33%  '
$ws.Range("E37").Value = '  +17.86%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E39").Value = '  +5.13%  '
$ws.Range("E40").Value = '  +14.14%  '
$ws.Range("D41").Value = '3.31'
$ws.Range("E41").Value = '  +14.54%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").Value = '3.006.08'
$ws.Range("E43").Value = '  +4.20%  '
$ws.Range("D44").Value = '2.64'
$ws.Range("E44").Value = '  +4.85%  '
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  +13.18%  '
$ws.Range("D46").Value = '3.34'
$ws.Range("E46").Value = '  +14.84%  '
$ws.Range("E47").Value = '  +7.73%  '
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("E49").Value = '  +6.96%  '
$ws.Range("D50").Value = '8.71'
$ws.Range("E50").Value = '  +15.87%  '
$ws.Range("D51").Value = '139.98'
$ws.Range("E51").Value = '  +5.21%  '
